$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 2081.889
$ws.Range("I115").Value = 2081.889
$ws.Range("K115").Value = 6245.667
$ws.Range("M115").Value = -4678.667

$ws.Range("H132").Value = 2718.6177
$ws.Range("I132").Value = 2782.1724
$ws.Range("K132").Value = 8346.5172
$ws.Range("M132").Value = -5816.5172

$ws.Range("H138").Value = 3311.7856
$ws.Range("I138").Value = 1210.381
$ws.Range("J138").Value = 4572.6284
$ws.Range("K138").Value = 3631.143
$ws.Range("L138").Value = 13717.8852
$ws.Range("M138").Value = 1508.857
$ws.Range("N138").Value = -23997.8852

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9291.833000000001
$ws.Range("I45").Value = 11535
$ws.Range("K45").Value = 11535
$ws.Range("M45").Value = -11158

$ws.Range("H74").Value = 1852.125
$ws.Range("I74").Value = 2128.3572
$ws.Range("J74").Value = 1465.4
$ws.Range("K74").Value = 2128.3572
$ws.Range("L74").Value = 1465.4
$ws.Range("M74").Value = -1254.3572
$ws.Range("N74").Value = -3213.4

$ws.Range("H77").Value = 1852.125
$ws.Range("I77").Value = 2128.3572
$ws.Range("J77").Value = 1465.4
$ws.Range("K77").Value = 10641.786
$ws.Range("L77").Value = 7327
$ws.Range("M77").Value = -6273.786
$ws.Range("N77").Value = -16063

$ws.Range("H102").Value = 1987.3704
$ws.Range("I102").Value = 2049.6956
$ws.Range("K102").Value = 2049.6956
$ws.Range("M102").Value = -427.6956

$ws.Range("H122").Value = 2150.5454
$ws.Range("I122").Value = 1805.0555
$ws.Range("J122").Value = 3705.25
$ws.Range("K122").Value = 5415.166499999999
$ws.Range("L122").Value = 11115.75
$ws.Range("M122").Value = -2965.166499999999
$ws.Range("N122").Value = -16015.75

$ws.Range("H132").Value = 2258.4546
$ws.Range("I132").Value = 2258.4546
$ws.Range("K132").Value = 6775.3638
$ws.Range("M132").Value = -4245.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1454
$ws.Range("I20").Value = 2051.6875
$ws.Range("J20").Value = 584.63635
$ws.Range("K20").Value = 2051.6875
$ws.Range("L20").Value = 584.63635
$ws.Range("M20").Value = -1804.6875
$ws.Range("N20").Value = -1078.63635

$ws.Range("H86").Value = 3586.5715
$ws.Range("I86").Value = 2841.182
$ws.Range("J86").Value = 6319.6665
$ws.Range("K86").Value = 2841.182
$ws.Range("L86").Value = 6319.6665
$ws.Range("M86").Value = -1718.182
$ws.Range("N86").Value = -8565.666499999999

$ws.Range("H89").Value = 3586.5715
$ws.Range("I89").Value = 2841.182
$ws.Range("J89").Value = 6319.6665
$ws.Range("K89").Value = 14205.91
$ws.Range("L89").Value = 31598.3325
$ws.Range("M89").Value = -8589.91
$ws.Range("N89").Value = -42830.3325

$ws.Range("H107").Value = 3125.05
$ws.Range("I107").Value = 2969.3572
$ws.Range("K107").Value = 2969.3572
$ws.Range("M107").Value = -1049.3572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 522.2222
$ws.Range("J22").Value = 166.66667
$ws.Range("L22").Value = 166.66667
$ws.Range("N22").Value = -866.6666700000001

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws.Range("H134").Value = 4343.436
$ws.Range("I134").Value = 3878.258
$ws.Range("J134").Value = 6146
$ws.Range("K134").Value = 11634.774
$ws.Range("L134").Value = 18438
$ws.Range("M134").Value = -9099.773999999999
$ws.Range("N134").Value = -23508

$ws.Range("H141").Value = 68391.625
$ws.Range("J141").Value = 68391.625
$ws.Range("L141").Value = 68391.625
$ws.Range("N141").Value = -78751.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 6325.636
$ws.Range("I118").Value = 3474.75
$ws.Range("J118").Value = 7954.7144
$ws.Range("K118").Value = 10424.25
$ws.Range("L118").Value = 23864.1432
$ws.Range("M118").Value = -9181.25
$ws.Range("N118").Value = -26350.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6049.7856
$ws.Range("I70").Value = 6103.6665
$ws.Range("J70").Value = 5726.5
$ws.Range("K70").Value = 6103.6665
$ws.Range("L70").Value = 5726.5
$ws.Range("M70").Value = -5833.6665
$ws.Range("N70").Value = -6266.5

$ws.Range("H73").Value = 6049.7856
$ws.Range("I73").Value = 6103.6665
$ws.Range("J73").Value = 5726.5
$ws.Range("K73").Value = 6103.6665
$ws.Range("L73").Value = 5726.5
$ws.Range("M73").Value = -5167.6665
$ws.Range("N73").Value = -7598.5

$ws.Range("H132").Value = 3137.0557
$ws.Range("I132").Value = 3302.6155
$ws.Range("K132").Value = 9907.8465
$ws.Range("M132").Value = -7377.8465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2132
$ws.Range("I46").Value = 1998
$ws.Range("J46").Value = 2400
$ws.Range("K46").Value = 1998
$ws.Range("L46").Value = 2400
$ws.Range("M46").Value = -1810
$ws.Range("N46").Value = -2776

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10950
$ws.Range("J41").Value = 10950
$ws.Range("L41").Value = 10950
$ws.Range("N41").Value = -11730

$ws.Range("H45").Value = 19835.285
$ws.Range("I45").Value = 16941
$ws.Range("J45").Value = 20993
$ws.Range("K45").Value = 16941
$ws.Range("L45").Value = 20993
$ws.Range("M45").Value = -16450
$ws.Range("N45").Value = -21975

$ws.Range("H62").Value = 18836.875
$ws.Range("I62").Value = 15199
$ws.Range("J62").Value = 24900
$ws.Range("K62").Value = 15199
$ws.Range("L62").Value = 24900
$ws.Range("M62").Value = -14575
$ws.Range("N62").Value = -26148

$ws.Range("H65").Value = 18836.875
$ws.Range("I65").Value = 15199
$ws.Range("J65").Value = 24900
$ws.Range("K65").Value = 75995
$ws.Range("L65").Value = 124500
$ws.Range("M65").Value = -72875
$ws.Range("N65").Value = -130740

$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30630

$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32184

$ws.Range("H100").Value = 4444.5454
$ws.Range("I100").Value = 6358
$ws.Range("K100").Value = 12716
$ws.Range("M100").Value = -12175

$ws.Range("H122").Value = 2592.5386
$ws.Range("I122").Value = 2551.3157
$ws.Range("J122").Value = 2704.4285
$ws.Range("K122").Value = 7653.9471
$ws.Range("L122").Value = 8113.2855
$ws.Range("M122").Value = -5203.9471
$ws.Range("N122").Value = -13013.2855

$ws.Range("H126").Value = 2906.5334
$ws.Range("I126").Value = 1719.6
$ws.Range("K126").Value = 5158.799999999999
$ws.Range("M126").Value = -2688.799999999999

$ws.Range("H136").Value = 2653.5908
$ws.Range("I136").Value = 1717.875
$ws.Range("J136").Value = 5148.8335
$ws.Range("K136").Value = 5153.625
$ws.Range("L136").Value = 15446.5005
$ws.Range("M136").Value = -2603.625
$ws.Range("N136").Value = -20546.5005
